# Major Change in Learning Rate Scheduler
# Adds a "File Name"/"Details" lookup block (H:I, plus a K helper column)
# above the existing data, and appends a new "pant" (pa1/pa2/pa3) group of
# 30 rows (121-150) to the existing A:E file-count table.
#
# NOTE: the order in which brand-new literal strings are first written
# determines their slot in the shared-strings table, so the writes below
# are sequenced to reproduce that exact table (towel/shirt/sweater/tshirt
# before the pant-group rows, "pant" itself only after them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Headers for the new summary table.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 8).Value  = "File Name"
$ws.Cells.Item(1, 9).Value  = "Details"
$ws.Cells.Item(1, 11).Value = "File Name"

# ---------------------------------------------------------------------
# 2) Fill in the H (File Name) / I (Details) rows for the four
#    pre-existing groups (rows 2-13; rows 14-16 = "pant" are filled
#    later). "towel" is written first so it mints the lowest shared-
#    string index of the four, matching the saved workbook.
# ---------------------------------------------------------------------
$hi = @(
  ,@(8,  "towel",   10)
  ,@(2,  "shirt",   10)
  ,@(3,  "shirt",   10)
  ,@(4,  "shirt",   10)
  ,@(5,  "sweater", 10)
  ,@(6,  "sweater",  9)
  ,@(7,  "sweater", 10)
  ,@(9,  "towel",   10)
  ,@(10, "towel",   10)
  ,@(11, "tshirt",  10)
  ,@(12, "tshirt",  10)
  ,@(13, "tshirt",  10)
)
foreach ($r in $hi) {
  $ws.Cells.Item($r[0], 8).Value = $r[1]
  $ws.Cells.Item($r[0], 9).Value = $r[2]
}

# ---------------------------------------------------------------------
# 3) K helper column (rows 2-5 reuse existing group names).
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 11).Value = "shirt"
$ws.Cells.Item(3, 11).Value = "sweater"
$ws.Cells.Item(4, 11).Value = "towel"
$ws.Cells.Item(5, 11).Value = "tshirt"

# ---------------------------------------------------------------------
# 4) Append the new "pant" group - rows 121-150 - to the A:E table.
#    This is where string indices 131-160 (pa1m1 ... pa3m10) get minted.
# ---------------------------------------------------------------------
$pant = @(
  ,@("pa1m1",  170, 158)
  ,@("pa1m2",   59,  68)
  ,@("pa1m3",  116, 154)
  ,@("pa1m4",  184, 151)
  ,@("pa1m5",  192, 163)
  ,@("pa1m6",  191, 155)
  ,@("pa1m7",   59,  63)
  ,@("pa1m8",  124, 135)
  ,@("pa1m9",   78,  83)
  ,@("pa1m10",  79,  81)
  ,@("pa2m1",   97,  94)
  ,@("pa2m2",  189, 192)
  ,@("pa2m3",  189, 217)
  ,@("pa2m4",  230, 211)
  ,@("pa2m5",   68,  71)
  ,@("pa2m6",  120, 115)
  ,@("pa2m7",   69,  58)
  ,@("pa2m8",  100, 101)
  ,@("pa2m9",  100, 102)
  ,@("pa2m10", 120, 120)
  ,@("pa3m1",  142, 165)
  ,@("pa3m2",  184, 218)
  ,@("pa3m3",  157, 141)
  ,@("pa3m4",  170, 143)
  ,@("pa3m5",  193, 218)
  ,@("pa3m6",  134, 117)
  ,@("pa3m7",  107, 114)
  ,@("pa3m8",   36,  36)
  ,@("pa3m9",   87,  86)
  ,@("pa3m10", 119, 122)
)

$row = 121
foreach ($r in $pant) {
  $ws.Cells.Item($row, 1).Value = $r[0]
  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $row++
}

# D121:D150 / E121:E150 carry on the same "B+C" / "D*3" pattern as the
# rest of the column - set once across the whole new block so they come
# out as one shared-formula group (matching D68:D120/E67:E120 above).
$ws.Range("D121:D150").Formula = "=B121+C121"
$ws.Range("E121:E150").Formula = "=D121*3"

# ---------------------------------------------------------------------
# 5) Finally, the "pant" rows of the summary table - this string is
#    minted last (index 161), after the pa1/pa2/pa3 detail rows above.
# ---------------------------------------------------------------------
$ws.Cells.Item(14, 8).Value = "pant"
$ws.Cells.Item(14, 9).Value = 10
$ws.Cells.Item(15, 8).Value = "pant"
$ws.Cells.Item(15, 9).Value = 10
$ws.Cells.Item(16, 8).Value = "pant"
$ws.Cells.Item(16, 9).Value = 10
$ws.Cells.Item(6, 11).Value = "pant"

# ---------------------------------------------------------------------
# 6) Column widths for the new H, I, K columns.
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth  = 18.7
$ws.Columns.Item(9).ColumnWidth  = 14.83
$ws.Columns.Item(11).ColumnWidth = 9.83

# ---------------------------------------------------------------------
# 7) Selection ends on B150, matching the author's final cursor position.
# ---------------------------------------------------------------------
$ws.Range("B150").Select()
